# Refactor synthetic array: insert a new "statut_name" column (French
# human-readable label for the "statut_label" column) right after
# "statut_label", shifting NCTId..intervention_type one column to the
# right (C:L -> D:M).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at C; this shifts C:L -> D:M automatically,
# including the dimension (A1:L16 -> A1:M16).
$ws.Columns("C:C").Insert()

# Header for the newly inserted column: copy the look of the
# neighbouring header cell, then set its text.
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C1").Value = "statut_name"

# Map each "statut_label" (column B) value to its French description.
$statusMap = @{
    "noir"   = "pas de résultat ni de publication"
    "vert"   = "résultat et / ou publication posté dans les 12 mois"
    "orange" = "résultat et / ou publication posté dans les 36 mois"
    "rouge"  = "résultat et / ou publication posté"
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, "B").End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $label = $ws.Cells.Item($r, 2).Text
    $ws.Cells.Item($r, 3).Value = $statusMap[$label]
}
